# Stage 6. A check for uniqueness is added. The error page was changed.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Splin")
$ws2 = $wb.Worksheets.Item("Iriao")

# --- Sheet "Splin": insert two rows before the old row 3 (Alexander Vasiliev /
# Georgia entry) so that the original data ends up on row 5, then fill rows
# 3 and 4 with duplicate "Alexander Vasiliev / Russia" entries used for the
# new uniqueness check, and patch row 5's Gender/Phone to the corrected values.
$ws1.Rows.Item(3).Insert() | Out-Null
$ws1.Rows.Item(3).Insert() | Out-Null

$ws1.Range("A3").Value = "Alexander Vasiliev"
$ws1.Range("B3").Value = 1
$ws1.Range("C3").Value = 25399
$ws1.Range("D3").Value = 7345934509
$ws1.Range("E3").Value = "Russia"
$ws1.Range("F3").Value = "Russian"
$ws1.Range("G3").Value = "Moscow"

$ws1.Range("A4").Value = "Alexander Vasiliev"
$ws1.Range("B4").Value = 1
$ws1.Range("C4").Value = 21746
$ws1.Range("D4").Value = 7345934509
$ws1.Range("E4").Value = "Russia"
$ws1.Range("F4").Value = "Russian"
$ws1.Range("G4").Value = "Moscow"

$ws1.Range("B5").Value = 1
$ws1.Range("D5").Value = 7345934509

# Stray customWidth column 8 definition left over from the author's edit.
$ws1.Columns.Item(8).ColumnWidth = 8.88671875

# --- Sheet "Iriao": append a new row 4 (duplicate-name check case, with a
# deliberately non-numeric Gender value "авва").
$ws2.Range("C3").Copy() | Out-Null
$ws2.Range("C4").PasteSpecial(-4122) | Out-Null

$ws2.Range("A4").Value = "Birdzina Muкia"
$ws2.Range("B4").Value = "авва"
$ws2.Range("C4").Value = 29906
$ws2.Range("D4").Value = 6756453423
$ws2.Range("E4").Value = "Georgia"
$ws2.Range("F4").Value = "Georgian"
$ws2.Range("G4").Value = "Tbilisi"

# --- Restore / update view selections on both sheets (Iriao stays the active tab).
$ws1.Range("F14").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("G8").Select() | Out-Null
